# "Change module to lecture"
# - Rename the "Modules" sheet to "Lectures".
# - Update the Cells!C2:C39 list-validation so it points at the renamed
#   sheet instead of the old name.
# - Make "Lectures" the active/selected sheet/tab (it was "Domains" before).

$wb = $excel.ActiveWorkbook

$modulesSheet = $wb.Worksheets.Item("Modules")
$modulesSheet.Name = "Lectures"

$cellsSheet = $wb.Worksheets.Item("Cells")
$cellsSheet.Range("C2:C39").Validation.Modify(3, 1, 1, "=Lectures!`$A`$2:`$A`$11")

$modulesSheet.Activate()
